# Add additional dimension item columns for b indicator
#
# Mirrors the existing "A.add_dim_1" / "A.add_dim_1_uid" / "A.add_dim_1_items" /
# "A.add_dim_1_items_uid" columns (M:P) but for the B indicator, inserting four
# new columns right before the trailing "calculation" column on the
# data_required sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_required")

# Insert 4 blank columns just before the old last column (V, "calculation"),
# which pushes "calculation" to Z.
$ws.Columns("V:Y").Insert()

# The engine copies the left-neighbour's formatting into freshly inserted
# columns; reset it back to the workbook default so the new cells are
# unstyled (matching the plain cells used elsewhere for these helper columns).
$ws.Range("V1:Y55").Style = "Normal"

# New header row (row 1) labels for the inserted columns.
$ws.Range("V1").Value = "B.add_dim_1"
$ws.Range("W1").Value = "B.add_dim_1_uid"
$ws.Range("X1").Value = "B.add_dim_1_items"
$ws.Range("Y1").Value = "B.add_dim_1_items_uid"

# Data rows (2-55): default every new cell to "NA", same as the other
# optional dimension columns.
$ws.Range("V2:Y55").Value = "NA"

# The worksheet's AutoFilter range needs to grow from A1:V55 to A1:Z55.
$ws.AutoFilterMode = $false
$ws.Range("A1:Z55").AutoFilter()

# Keep the hidden _FilterDatabase defined name for this sheet in sync with
# the new autofilter range.
$wb.Names.Item("data_required!_FilterDatabase").RefersTo = "=data_required!`$A`$1:`$Z`$55"

# Approximate the editor's final on-screen selection.
$ws.Range("Y61").Select()
